$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 (even_MAG-GUT449.fa) entirely; subsequent rows shift up.
$ws.Rows.Item(2).Delete()
